$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-ambiguous string updates
$ws.Range("D2").Value = "63.667.36"
$ws.Range("E2").Value = "  -1.91%  "
$ws.Range("D3").Value = "3.043.20"
$ws.Range("E3").Value = "  -2.09%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("E6").Value = "  -2.01%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.044.95"
$ws.Range("E8").Value = "  -1.87%  "
$ws.Range("E9").Value = "  +3.81%  "
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("E11").Value = "  -12.56%  "
$ws.Range("E12").Value = "  +4.61%  "
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").Value = "3.540.20"
$ws.Range("E15").Value = "  -1.79%  "
$ws.Range("D16").Value = "63.720.37"
$ws.Range("E16").Value = "  -1.73%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.030.89"
$ws.Range("E18").Value = "  -2.38%  "
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("E21").Value = "  +1.27%  "
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("E23").Value = "  +8.99%  "
$ws.Range("E24").Value = "  -1.12%  "
$ws.Range("E25").Value = "  +1.77%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("E28").Value = "  -1.61%  "
$ws.Range("E29").Value = "  -2.53%  "
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("E32").Value = "  -1.62%  "
$ws.Range("E33").Value = "  -1.58%  "
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("E35").Value = "  -1.02%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  -0.80%  "
$ws.Range("E38").Value = "  -6.56%  "
$ws.Range("E39").Value = "  -2.24%  "
$ws.Range("D40").Value = "2.997.04"
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("E41").Value = "  -4.39%  "
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("E44").Value = "  +2.52%  "
$ws.Range("E45").Value = "  -2.34%  "
$ws.Range("E46").Value = "  +5.72%  "
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("D50").Value = "0.0₃0510"
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("E51").Value = "  -0.58%  "

# Numeric-looking text values: force text format so Excel keeps them as strings
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "555.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.519"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.152"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.23"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.486"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000228"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.109"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "473.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.681"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0407"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "440.47"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0810"
$ws.Range("D39").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.268"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.08"
$ws.Range("D51").Style = "Normal"

Write-Output "Applied cryptos update"
